# Add a new, blank slide at the end of the deck (slide 18).
# This mirrors inserting a new "Leer" (Blank) layout slide after the
# current last slide, via Slides.Add with the ppLayoutBlank layout.

$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 12)  # 12 = ppLayoutBlank
